{"js": "// Office.js (Word JavaScript API) \u2014 apply the IPSSI -> EFREI /\n// \"DEV,DATA,IA\" -> \"R\u00e9seaux et S\u00e9curit\u00e9\" program-name change described\n// by the diff.\n//\n// Strategy: locate each affected phrase with Range.search() (which\n// matches across run boundaries) and rewrite it in place with\n// Range.insertText(text, \"Replace\"). Every needle below is unique in\n// the document, so each search returns exactly one hit.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${searchText}`);\n  }\n\n  // Replace the (first / only) match; keep the formatting of its run.\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1) Objet line: \"Apprentissage parcours DEV,DATA,IA - IPSSI \"\n//    -> \"Apprentissage parcours R\u00e9seaux et S\u00e9curit\u00e9 - EFREI\"\nawait replaceOnce(\"DEV,DATA,IA - IPSSI \", \"R\u00e9seaux et S\u00e9curit\u00e9 - EFREI\");\n\n// 2) \"IPSSI est une \u00e9cole reconnue...\" -> \"EFREI est une \u00e9cole reconnue...\"\nawait replaceOnce(\"IPSSI est\", \"EFREI est\");\n\n// 3) \"...votre programme parcours DEV, DATA,IA en alternance.\"\n//    -> \"...votre programme parcours R\u00e9seaux et S\u00e9curit\u00e9 en alternance.\"\nawait replaceOnce(\"DEV, DATA,IA \", \"R\u00e9seaux et S\u00e9curit\u00e9 \");\n\n// 4) \"...un grand int\u00e9r\u00eat pour le Web.\" -> \"...pour le R\u00e9seaux et la S\u00e9curit\u00e9.\"\nawait replaceOnce(\"pour le Web.\", \"pour le R\u00e9seaux et la S\u00e9curit\u00e9.\");\n\n// 5) \"...et sur des scripts en Python, JAVA mais aussi sur le cloud, j'ai...\"\n//    -> \"...et sur des scripts en Python mais aussi dans le R\u00e9seaux et la\n//    S\u00e9curit\u00e9, j'ai...\"\nawait replaceOnce(\n  \"Python, JAVA mais aussi sur le cloud, \",\n  \"Python mais aussi dans le R\u00e9seaux et la S\u00e9curit\u00e9, \"\n);\n\n// 6) \"...int\u00e9grer une \u00e9cole comme IPSSI.\" -> \"...comme EFREI.\"\nawait replaceOnce(\"comme IPSSI.\", \"comme EFREI.\");\n", "ps1": "# Word COM interop (PowerShell-style) \u2014 apply the IPSSI -> EFREI /\n# \"DEV,DATA,IA\" -> \"R\u00e9seaux et S\u00e9curit\u00e9\" program-name change described\n# by the diff.\n#\n# Strategy: use Range.Find.Execute(..., Replace:=wdReplaceAll) against a\n# fresh $d.Content range for each phrase. Every FindText below is unique\n# in the document, so each call rewrites exactly one spot, and\n# Find.Execute happily matches text that spans multiple runs.\n\n$d = $word.ActiveDocument\n\n# Positional Find.Execute signature used throughout:\n#   (FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#    MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll).\n\n# 1) Objet line: \"Apprentissage parcours DEV,DATA,IA - IPSSI \"\n#    -> \"Apprentissage parcours R\u00e9seaux et S\u00e9curit\u00e9 - EFREI\"\n$d.Content.Find.Execute(\"DEV,DATA,IA - IPSSI \", $true, $false, $false, $false, $false, $true, 1, $false, \"R\u00e9seaux et S\u00e9curit\u00e9 - EFREI\", 2)\n\n# 2) \"IPSSI est une \u00e9cole reconnue...\" -> \"EFREI est une \u00e9cole reconnue...\"\n$d.Content.Find.Execute(\"IPSSI est\", $true, $false, $false, $false, $false, $true, 1, $false, \"EFREI est\", 2)\n\n# 3) \"...votre programme parcours DEV, DATA,IA en alternance.\"\n#    -> \"...votre programme parcours R\u00e9seaux et S\u00e9curit\u00e9 en alternance.\"\n$d.Content.Find.Execute(\"DEV, DATA,IA \", $true, $false, $false, $false, $false, $true, 1, $false, \"R\u00e9seaux et S\u00e9curit\u00e9 \", 2)\n\n# 4) \"...un grand int\u00e9r\u00eat pour le Web.\" -> \"...pour le R\u00e9seaux et la S\u00e9curit\u00e9.\"\n$d.Content.Find.Execute(\"pour le Web.\", $true, $false, $false, $false, $false, $true, 1, $false, \"pour le R\u00e9seaux et la S\u00e9curit\u00e9.\", 2)\n\n# 5) \"...et sur des scripts en Python, JAVA mais aussi sur le cloud, j'ai...\"\n#    -> \"...et sur des scripts en Python mais aussi dans le R\u00e9seaux et la\n#    S\u00e9curit\u00e9, j'ai...\"\n$d.Content.Find.Execute(\"Python, JAVA mais aussi sur le cloud, \", $true, $false, $false, $false, $false, $true, 1, $false, \"Python mais aussi dans le R\u00e9seaux et la S\u00e9curit\u00e9, \", 2)\n\n# 6) \"...int\u00e9grer une \u00e9cole comme IPSSI.\" -> \"...comme EFREI.\"\n$d.Content.Find.Execute(\"comme IPSSI.\", $true, $false, $false, $false, $false, $true, 1, $false, \"comme EFREI.\", 2)\n"}
